# SS-OBS-DB.xlsx update: "Major Update - Added 2021"
# Adds 2021 Mars observation rows (118-133) to the Mars sheet, reusing the
# same column layout / formulas as the existing 2020 rows (114-117) above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mars")

# --- Row-by-row source data transcribed from the updated report --------
$rowNums = @(118, 119, 120, 121, 122, 123, 124, 125, 126, 127, 128, 129, 130, 131, 132, 133)

$C = @(44076.25, 44078.5, 44080.25, 44090.229166666664, 44099.208333333336, 44109.208333333336, 44132.104166666664, 44132.229166666664, 44133.145833333336, 44137.145833333336, 44138.166666666664, 44141.15625, 44152.125, 44153.21875, 44164.114583333336, 44172.145833333336)
$D = @(0.08, 0.07, 0.07, 0.04, 0.02, 0.01, 0.01, 0.01, 0.01, 0.02, 0.02, 0.03, 0.05, 0.05, 0.07, 0.09)
$E = @(275.4, 344.7, 238.7, 140.5, 52.6, 324.1, 84.7, 128.5, 90.4, 54.6, 53, 22.3, 271.3, 295, 156.6, 92.5)
$F = @(270, 272, 273, 279, 284, 291, 304, 304, 305, 307, 308, 310, 316, 317, 323, 327)
$G = @(19.1, 19.4, 19.7, 21.1, 22.1, 22.6, 20.7, 20.7, 20.5, 19.9, 19.7, 19.1, 17.1, 16.9, 14.9, 13.6)
$H = @(-1.8, -1.9, -1.9, -2.2, -2.4, -2.4, -2.3, -2.3, -2.2, -2.1, -2.1, -2.0, -1.6, -1.6, -1.2, -0.9)
$M = @(27, 46, 30, 32, 33, 41, 35, 55, 45, 48, 52, 52, 52, 52, 55, 57)
$N = @(0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8, 0.8)

for ($i = 0; $i -lt $rowNums.Length; $i++) {
    $r = $rowNums[$i]

    # Raw inputs
    $ws.Cells.Item($r, 3).Value  = $C[$i]     # C - date/time of observation
    $ws.Cells.Item($r, 4).Value  = $D[$i]     # D - phase / fraction
    $ws.Cells.Item($r, 5).Value  = $E[$i]     # E - Central Meridian longitude
    $ws.Cells.Item($r, 6).Value  = $F[$i]     # F - Ls (areocentric longitude)
    $ws.Cells.Item($r, 7).Value  = $G[$i]     # G - apparent diameter
    $ws.Cells.Item($r, 8).Value  = $H[$i]     # H - magnitude
    $ws.Cells.Item($r, 10).Value = "X"        # J - imaged flag (shared string "X")
    $ws.Cells.Item($r, 13).Value = $M[$i]     # M - altitude
    $ws.Cells.Item($r, 15).Value = 0          # O
    $ws.Cells.Item($r, 16).Value = 1          # P
    $ws.Cells.Item($r, 18).Value = 44053      # R - reference opposition date
    $ws.Cells.Item($r, 19).Value = 0.42       # S - aperture (m)
    $ws.Cells.Item($r, 22).Value = 22.22      # V - site latitude

    # N - seeing estimate: same value as other columns, but flagged in red
    # (new cellXf: numFmtId 12, red font, centered) to mark rows added in
    # this update.
    $nCell = $ws.Cells.Item($r, 14)
    $nCell.Value = $N[$i]
    $nCell.NumberFormat = "# ?/?"
    $nCell.Font.Name = "Arial"
    $nCell.Font.Size = 10
    $nCell.Font.Color = 255
    $nCell.HorizontalAlignment = -4108
    $nCell.VerticalAlignment = -4160

    # Formulas (identical shape to the pre-existing rows 114-117)
    $ws.Cells.Item($r, 9).Formula  = "=-LOG((1/(G$r^2))*(2.511^(-H$r)))/LOG(2.511)"
    $ws.Cells.Item($r, 20).Formula = "=1/SIN(RADIANS(M$r+244/(165+47*M$r^1.1)))"
    $ws.Cells.Item($r, 21).Formula = "=DEGREES(ASIN(SIN(RADIANS(F$r))*SIN(RADIANS(25.19))))"
    $ws.Cells.Item($r, 23).Formula = "=668.5921*MOD((C$r-Z`$1),686.9726)/686.9726"
    $ws.Cells.Item($r, 24).Formula = "=(S$r/(G$r/2))*T$r/(IF(ISBLANK(N$r),0.5,N$r))"
    $ws.Cells.Item($r, 25).Formula = "=0.1/X$r"
}
